$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.644.39"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "2.124.62"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.86"
$ws.Range("E5").Value = "  +5.07%  "

$ws.Range("E6").Value = "  +0.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5274"
$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4531"
$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.03"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09125"
$ws.Range("E10").Value = "  +2.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.185"
$ws.Range("E11").Value = "  +0.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.63"
$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").Value = "2.125.09"
$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.862"
$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.106"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.45"

$ws.Range("E17").Value = "  +3.23%  "

$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06721"
$ws.Range("E19").Value = "  +1.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.52"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.347"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").Value = "30.712.49"
$ws.Range("E23").Value = "  +0.67%  "

$ws.Range("E24").Value = "  +3.05%  "

$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").Value = "2.363.04"
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.52"
$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.578"
$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.32"
$ws.Range("E29").Value = "  +0.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.49"
$ws.Range("E30").Value = "  +2.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.201"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1083"
$ws.Range("E32").Value = "  +0.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.666"
$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.402"
$ws.Range("E34").Value = "  +1.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.031"
$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.38"
$ws.Range("E36").Value = "  -1.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.980"
$ws.Range("E37").Value = "  +5.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02662"
$ws.Range("E38").Value = "  +2.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06894"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("E40").Value = "  +0.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.59"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6942"
$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.280"
$ws.Range("E43").Value = "  +2.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.88"
$ws.Range("E44").Value = "  +5.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6491"
$ws.Range("E45").Value = "  +1.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.339"
$ws.Range("E46").Value = "  -0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.764"
$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000367"
$ws.Range("E48").Value = "  +7.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.258"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07321"
$ws.Range("E50").Value = "  +2.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.93"
$ws.Range("E51").Value = "  -0.57%  "
